$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")

# Revert the previously-added hyperlink on B6: remove the hyperlink object
# (relationship) and restore the old literal URL text that used to live there
# before the hyperlink was introduced.
$wsAbout.Hyperlinks.Delete()
$wsAbout.Range("B6").Value = "http://yosemite.epa.gov/EE%5Cepa%5Ceed.nsf/webpages/MortalityRiskValuation.html#whatvalue"
